$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.497.27'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '1.887.86'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.82'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.79%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4691'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2900'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06496'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.31'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07749'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("D12").Value = '1.888.29'
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.75'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7281'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.192'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '282.20'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.29%  '
$ws.Range("D17").Value = '30.499.02'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.04'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.29%  '
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007490'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("D21").Value = '2.132.72'
$ws.Range("E21").Value = '  +0.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.273'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.272'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '163.63'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.092'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.95'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.896'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.334'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09717'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.471'
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.287'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.120'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04864'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.127'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6946'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.715'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01898'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.848'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.77'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.211'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.001'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4259'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8267'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.47'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.570'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.979'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.17'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '916.72'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05757'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.86%  '
